# Auto-generated Excel COM-interop script
# Applies updated market-price / profit figures to the Leve profit sheets
# (mirrors a scheduled data-refresh run across all 8 sheets).

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 226.66667
$ws.Range("I2").Value = 232
$ws.Range("J2").Value = 200
$ws.Range("K2").Value = 232
$ws.Range("L2").Value = 200
$ws.Range("M2").Value = -119
$ws.Range("N2").Value = -426
$ws.Range("H4").Value = 115.333336
$ws.Range("I4").Value = 115.333336
$ws.Range("K4").Value = 115.333336
$ws.Range("M4").Value = -1.333336000000003
$ws.Range("H18").Value = 650
$ws.Range("I18").Value = 650
$ws.Range("J18").Value = 0
$ws.Range("K18").Value = 650
$ws.Range("L18").Value = 0
$ws.Range("M18").Value = -366
$ws.Range("N18").ClearContents()
$ws.Range("H62").Value = 31257666
$ws.Range("I62").Value = 1100
$ws.Range("K62").Value = 1100
$ws.Range("M62").Value = -476
$ws.Range("H65").Value = 31257666
$ws.Range("I65").Value = 1100
$ws.Range("K65").Value = 5500
$ws.Range("M65").Value = -2380
$ws.Range("H96").Value = 305.27585
$ws.Range("I96").Value = 249.88889
$ws.Range("J96").Value = 395.9091
$ws.Range("K96").Value = 749.6666700000001
$ws.Range("L96").Value = 1187.7273
$ws.Range("M96").Value = 623.3333299999999
$ws.Range("N96").Value = -3933.7273
$ws.Range("H129").Value = 731.1111
$ws.Range("I129").Value = 408.54544
$ws.Range("J129").Value = 952.875
$ws.Range("K129").Value = 1225.63632
$ws.Range("L129").Value = 2858.625
$ws.Range("M129").Value = 3774.36368
$ws.Range("N129").Value = -12858.625
$ws.Range("H137").Value = 1952.0222
$ws.Range("I137").Value = 1190.125
$ws.Range("J137").Value = 3827.4614
$ws.Range("K137").Value = 3570.375
$ws.Range("L137").Value = 11482.3842
$ws.Range("M137").Value = -1020.375
$ws.Range("N137").Value = -16582.3842

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H80").Value = 37950
$ws.Range("J80").Value = 37950
$ws.Range("L80").Value = 37950
$ws.Range("N80").Value = -39946
$ws.Range("H83").Value = 37950
$ws.Range("J83").Value = 37950
$ws.Range("L83").Value = 113850
$ws.Range("N83").Value = -123834
$ws.Range("H122").Value = 1473.8572
$ws.Range("I122").Value = 1248.5454
$ws.Range("J122").Value = 2300
$ws.Range("K122").Value = 3745.6362
$ws.Range("L122").Value = 6900
$ws.Range("M122").Value = -1295.6362
$ws.Range("N122").Value = -11800
$ws.Range("H123").Value = 56460
$ws.Range("J123").Value = 56460
$ws.Range("L123").Value = 56460
$ws.Range("N123").Value = -66260

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H81").Value = 18850
$ws.Range("J81").Value = 18850
$ws.Range("L81").Value = 18850
$ws.Range("N81").Value = -20972
$ws.Range("H82").Value = 64845
$ws.Range("I82").Value = 107103.336
$ws.Range("J82").Value = 33151.25
$ws.Range("K82").Value = 107103.336
$ws.Range("L82").Value = 33151.25
$ws.Range("M82").Value = -106720.336
$ws.Range("N82").Value = -33917.25
$ws.Range("H84").Value = 18850
$ws.Range("J84").Value = 18850
$ws.Range("L84").Value = 56550
$ws.Range("N84").Value = -67158
$ws.Range("H85").Value = 64845
$ws.Range("I85").Value = 107103.336
$ws.Range("J85").Value = 33151.25
$ws.Range("K85").Value = 107103.336
$ws.Range("L85").Value = 33151.25
$ws.Range("M85").Value = -105777.336
$ws.Range("N85").Value = -35803.25

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H41").Value = 14203
$ws.Range("J41").Value = 21038.334
$ws.Range("L41").Value = 21038.334
$ws.Range("N41").Value = -21894.334
$ws.Range("H51").Value = 9537.4
$ws.Range("J51").Value = 9537.4
$ws.Range("L51").Value = 9537.4
$ws.Range("N51").Value = -11009.4
$ws.Range("H61").Value = 9537.4
$ws.Range("J61").Value = 9537.4
$ws.Range("L61").Value = 9537.4
$ws.Range("N61").Value = -10233.4
$ws.Range("H88").Value = 29189.182
$ws.Range("J88").Value = 30624.6
$ws.Range("L88").Value = 30624.6
$ws.Range("N88").Value = -31436.6
$ws.Range("H91").Value = 29189.182
$ws.Range("J91").Value = 30624.6
$ws.Range("L91").Value = 30624.6
$ws.Range("N91").Value = -33432.6
$ws.Range("H107").Value = 891.65515
$ws.Range("I107").Value = 870.7143
$ws.Range("J107").Value = 946.625
$ws.Range("K107").Value = 870.7143
$ws.Range("L107").Value = 946.625
$ws.Range("M107").Value = 1049.2857
$ws.Range("N107").Value = -4786.625
$ws.Range("H109").Value = 11500
$ws.Range("J109").Value = 11500
$ws.Range("L109").Value = 11500
$ws.Range("N109").Value = -13580
$ws.Range("H134").Value = 5683.3716
$ws.Range("I134").Value = 6482.931
$ws.Range("J134").Value = 1818.8334
$ws.Range("K134").Value = 19448.793
$ws.Range("L134").Value = 5456.5002
$ws.Range("M134").Value = -16913.793
$ws.Range("N134").Value = -10526.5002

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 130.73529
$ws.Range("I4").Value = 98.5
$ws.Range("J4").Value = 372.5
$ws.Range("K4").Value = 295.5
$ws.Range("L4").Value = 1117.5
$ws.Range("M4").Value = -183.5
$ws.Range("N4").Value = -1341.5
$ws.Range("H5").Value = 410.10526
$ws.Range("I5").Value = 383.13333
$ws.Range("K5").Value = 1149.39999
$ws.Range("M5").Value = -1037.39999
$ws.Range("H40").Value = 96
$ws.Range("I40").Value = 96
$ws.Range("J40").Value = 0
$ws.Range("K40").Value = 384
$ws.Range("L40").Value = 0
$ws.Range("M40").Value = -315
$ws.Range("N40").ClearContents()
$ws.Range("H135").Value = 410.10526
$ws.Range("I135").Value = 383.13333
$ws.Range("K135").Value = 3448.19997
$ws.Range("M135").Value = -913.1999700000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 1875.6471
$ws.Range("I122").Value = 1898.9231
$ws.Range("J122").Value = 1800
$ws.Range("K122").Value = 5696.7693
$ws.Range("L122").Value = 5400
$ws.Range("M122").Value = -3246.7693
$ws.Range("N122").Value = -10300
$ws.Range("H123").Value = 34487
$ws.Range("J123").Value = 34487
$ws.Range("L123").Value = 34487
$ws.Range("N123").Value = -39387

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 11113766
$ws.Range("I7").Value = 20001978
$ws.Range("K7").Value = 20001978
$ws.Range("M7").Value = -20001866
$ws.Range("H46").Value = 1466.0526
$ws.Range("I46").Value = 1409.2593
$ws.Range("J46").Value = 1605.4546
$ws.Range("K46").Value = 1409.2593
$ws.Range("L46").Value = 1605.4546
$ws.Range("M46").Value = -1221.2593
$ws.Range("N46").Value = -1981.4546
$ws.Range("H126").Value = 11113766
$ws.Range("I126").Value = 20001978
$ws.Range("K126").Value = 60005934
$ws.Range("M126").Value = -60003464
$ws.Range("H128").Value = 0
$ws.Range("J128").Value = 0
$ws.Range("L128").Value = 0
$ws.Range("N128").ClearContents()
$ws.Range("H136").Value = 2036.6786
$ws.Range("I136").Value = 1681.08
$ws.Range("K136").Value = 5043.24
$ws.Range("M136").Value = -2493.24

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H109").Value = 27088.5
$ws.Range("J109").Value = 27088.5
$ws.Range("L109").Value = 27088.5
$ws.Range("N109").Value = -29862.5
